{"js": "// \"started doing practical results evaluation\"\n//\n// 1) Table 3 caption: \"Lentel\u0117 3.\" -> \"Lentel\u0117 nr. 3.\"\n// 2) Fill in the previously-empty \"practical speedup\" results for the\n//    N=48000/p=1 and N=96000/p=1 rows of the 4th table (practical\n//    speedups table: N dydis | Naudot\u0173 procesori\u0173 kiekis |\n//    Nuosekliosios dalies pagreit\u0117jimas | Lygiagre\u010diosios dalies\n//    pagreit\u0117jimas | Bendras algoritmo pagreit\u0117jimas).\n\n// --- 1. Fix the Table 3 caption -------------------------------------------------\nconst captionResults = context.document.body.search(\"Lentel\u0117 3.\", { matchCase: true });\ncaptionResults.load(\"text\");\nawait context.sync();\n\nif (captionResults.items.length > 0) {\n  captionResults.items[0].insertText(\"Lentel\u0117 nr. 3.\", Word.InsertLocation.replace);\n}\n\n// --- 2. Fill in the practical speedup values -------------------------------------\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\n// 4th table on the page (0-indexed -> index 3) is the \"Praktiniai\n// pagreit\u0117jimai\" results table.\nconst resultsTable = tables.items[3];\n\n// Row index 2 == N dydis 48000, 1 processor.\nresultsTable.getCell(2, 2).value = \"0,96\";\nresultsTable.getCell(2, 3).value = \"0,88\";\nresultsTable.getCell(2, 4).value = \"0,89\";\n\n// Row index 3 == N dydis 96000, 1 processor.\nresultsTable.getCell(3, 2).value = \"0,96\";\nresultsTable.getCell(3, 3).value = \"0,90\";\nresultsTable.getCell(3, 4).value = \"0,90\";\n\nawait context.sync();\n", "ps1": "# \"started doing practical results evaluation\"\n#\n# 1) Table 3 caption: \"Lentel\u0117 3.\" -> \"Lentel\u0117 nr. 3.\"\n# 2) Fill in the previously-empty \"practical speedup\" results for the\n#    N=48000/p=1 and N=96000/p=1 rows of the 4th table (practical\n#    speedups table: N dydis | Naudot\u0173 procesori\u0173 kiekis |\n#    Nuosekliosios dalies pagreit\u0117jimas | Lygiagre\u010diosios dalies\n#    pagreit\u0117jimas | Bendras algoritmo pagreit\u0117jimas).\n\n$d = $word.ActiveDocument\n\n# --- 1. Fix the Table 3 caption ------------------------------------------------\n# The caption paragraph starts with the literal run \"Lentel\u0117 \" followed\n# immediately by a SEQ field that renders the table number (\"3\"). We only\n# want to insert \"nr. \" right after \"Lentel\u0117 \" (leaving the field alone),\n# so find the paragraph by its rendered text and then edit a sub-range\n# instead of doing a blind Find/Replace (which would clash with the field).\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"Lentel\u0117 3.*\") {\n        $prefix = $p.Range.Duplicate\n        $prefix.SetRange($p.Range.Start, $p.Range.Start + 8)\n        if ($prefix.Text -eq \"Lentel\u0117 \") {\n            $prefix.Collapse(0)\n            $prefix.InsertBefore(\"nr. \")\n        }\n        break\n    }\n}\n\n# --- 2. Fill in the practical speedup values -----------------------------------\n# 4th table on the page is the \"Praktiniai pagreit\u0117jimai\" results table.\n$t = $d.Tables.Item(4)\n\n# Row 3 == N dydis 48000, 1 processor.\n$t.Cell(3, 3).Range.Text = \"0,96\"\n$t.Cell(3, 4).Range.Text = \"0,88\"\n$t.Cell(3, 5).Range.Text = \"0,89\"\n\n# Row 4 == N dydis 96000, 1 processor.\n$t.Cell(4, 3).Range.Text = \"0,96\"\n$t.Cell(4, 4).Range.Text = \"0,90\"\n$t.Cell(4, 5).Range.Text = \"0,90\"\n"}
